$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The trailing space in the header "ID " was a typo; correct it to "ID".
$ws.Range("A1").Value = "ID"
